# Apply the "Updated cryptos list" price/volume refresh (coinranking scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells that are plain decimals (e.g. "121.96") would be
# auto-detected by Excel as numbers; prefixing with an apostrophe keeps them
# as text, matching the sheet's existing inlineStr/shared-string storage.

$ws.Range("D2").Value = "43.649.83"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "2.273.17"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'121.96"
$ws.Range("E5").Value = "  +5.29%  "

$ws.Range("D6").Value = "'264.96"
$ws.Range("E6").Value = "  -1.64%  "

$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = "  +1.69%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").Value = "'0.621"
$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("D10").Value = "'48.04"
$ws.Range("E10").Value = "  -1.30%  "

$ws.Range("D11").Value = "'0.0941"
$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").Value = "'8.98"
$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("D14").Value = "'15.41"
$ws.Range("E14").Value = "  -2.78%  "

$ws.Range("D15").Value = "'0.892"
$ws.Range("E15").Value = "  +4.07%  "

$ws.Range("D16").Value = "2.617.81"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("D17").Value = "2.273.90"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").Value = "43.591.88"
$ws.Range("E18").Value = "  -0.22%  "

$ws.Range("D19").Value = "'0.0000109"
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("D20").Value = "'6.96"
$ws.Range("E20").Value = "  -1.30%  "

$ws.Range("D21").Value = "'72.24"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "'2.43"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").Value = "'235.00"
$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("D24").Value = "'9.49"
$ws.Range("E24").Value = "  -4.14%  "

$ws.Range("D25").Value = "'2.85"
$ws.Range("E25").Value = "  -3.94%  "

$ws.Range("E26").Value = "  +1.71%  "

$ws.Range("D27").Value = "'11.77"
$ws.Range("E27").Value = "  +1.65%  "

$ws.Range("D28").Value = "'41.90"
$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D31").Value = "'171.57"
$ws.Range("E31").Value = "  -2.26%  "

$ws.Range("D32").Value = "'21.61"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").Value = "'0.0907"
$ws.Range("E33").Value = "  -2.93%  "

$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("E35").Value = "  +1.73%  "

$ws.Range("D36").Value = "'0.0378"
$ws.Range("E36").Value = "  +4.12%  "

$ws.Range("D37").Value = "'4.66"
$ws.Range("E37").Value = "  -2.58%  "

$ws.Range("D38").Value = "'4.07"
$ws.Range("E38").Value = "  +5.36%  "

$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("E40").Value = "  +5.15%  "

$ws.Range("D41").Value = "'75.14"
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").Value = "'13.87"
$ws.Range("E42").Value = "  -2.42%  "

$ws.Range("D43").Value = "'0.237"
$ws.Range("E43").Value = "  -2.35%  "

$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").Value = "'1.36"
$ws.Range("E45").Value = "  -2.65%  "

$ws.Range("D46").Value = "'5.75"
$ws.Range("E46").Value = "  -9.54%  "

$ws.Range("D47").Value = "'73.70"
$ws.Range("E47").Value = "  +37.43%  "

$ws.Range("D48").Value = "'8.53"
$ws.Range("E48").Value = "  -3.15%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  +1.05%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.26"
$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("D51").Value = "'101.19"
$ws.Range("E51").Value = "  -1.06%  "
